$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: columns B:F, rows 2:25
$block1 = @(
@(1.02,1.028333286578941,1.036506160765283,1.028308089439665,1.045316223443975),
@(1.02,1.029264745061004,1.037196341351081,1.029098810820422,1.046146191492865),
@(1.02,1.029868095796984,1.037643300497118,1.029611384564347,1.046683933459939),
@(1.02,1.030121895208278,1.037831288166892,1.029827090249436,1.046910165357188),
@(1.02,1.030164518051776,1.037862857122597,1.029863321020171,1.046948160324018),
@(1.02,1.02987148648626,1.037645812061496,1.029614265973019,1.046686955734765),
@(1.02,1.028647945160686,1.03673933374085,1.028575125147773,1.04559656964865),
@(1.02,1.026496824634731,1.035144878708707,1.026751172832529,1.043680585440731),
@(1.02,1.025066120818483,1.034083937938518,1.025540102070436,1.042407005970272),
@(1.02,1.024447425944714,1.033625040336938,1.025016876093686,1.04185644135889),
@(1.02,1.024217737869335,1.033454661584798,1.024822704860063,1.041652074644083),
@(1.02,1.024267001167957,1.033491204929086,1.024864347160224,1.041695905748493),
@(1.02,1.024428437346015,1.033610955211318,1.025000822194873,1.041839545527649),
@(1.02,1.024527919803406,1.033684747464629,1.025084932626688,1.041928065007398),
@(1.02,1.025107198709706,1.034114404065152,1.02557485175834,1.042443564317493),
@(1.02,1.025470782556309,1.034384050436875,1.025882480945101,1.04276716658003),
@(1.02,1.025682932872436,1.034541378424488,1.026062029174008,1.042956005383516),
@(1.02,1.025755283844468,1.034595031264577,1.026123269655155,1.043020409253057),
@(1.02,1.0254317653732,1.034355114975053,1.025849463511307,1.042732438100388),
@(1.02,1.024380895017501,1.033575689620052,1.024960628767241,1.041797243383411),
@(1.02,1.023720881875455,1.033086075919961,1.024402814815714,1.041210045708527),
@(1.02,1.024070699570384,1.033345586974014,1.024698424228092,1.041521254244514),
@(1.02,1.025449395317832,1.034368189517166,1.025864382318433,1.042748130136606),
@(1.02,1.027052350972483,1.035556732897641,1.027221851787699,1.044175260987338)
)

# Block 2: columns I:M, rows 2:25
$block2 = @(
@(1.033753740601323,1.033486186819577,1.039299707786021,1.031125284585218,1.048084805741589),
@(1.033896067206445,1.034058302321531,1.039799820788631,1.031723960519061,1.048726141269683),
@(1.033986983326679,1.034428524607848,1.040123116547517,1.032111613505977,1.049141176486108),
@(1.034024921664882,1.034584170967199,1.040258954475734,1.032274646155722,1.049315667406149),
@(1.034031275077548,1.034610304930425,1.040281757792864,1.032302023733622,1.049344965750858),
@(1.033987491372582,1.034430604344264,1.040124931919165,1.032113791708887,1.049143508004379),
@(1.033802084234763,1.033679530003825,1.039468786945269,1.031327553361859,1.048301537233928),
@(1.033466372774382,1.032356285498189,1.038310248527629,1.029944219626492,1.046818310057032),
@(1.033236550216357,1.031474354797172,1.037536396374288,1.029023493683477,1.045829868174029),
@(1.033135615480893,1.031092538374197,1.037200970650404,1.028625179244212,1.045401969031022),
@(1.033097911072045,1.030950725705319,1.037076328313773,1.028477283713682,1.045243045094109),
@(1.03310600841377,1.030981144490532,1.037103066774633,1.028509005245097,1.045277134032359),
@(1.033132503159379,1.031080815864831,1.037190668694958,1.028612953000355,1.0453888319767),
@(1.033148799272371,1.031142228172945,1.037244636501229,1.028677006124454,1.045457655043057),
@(1.03324321905404,1.031499696169545,1.037558650357455,1.029049936314351,1.045858268716641),
@(1.033302066318483,1.031723944701691,1.037755532059625,1.02928396440041,1.046109591595862),
@(1.033336253878043,1.031854751310112,1.037870336686478,1.029420504253572,1.04625619379419),
@(1.033347887687591,1.0318993540063,1.037909476440289,1.029467066753982,1.046306182976704),
@(1.033295766735809,1.031699884303541,1.037734411935386,1.029258851748887,1.046082626003916),
@(1.033124706985447,1.031051464817989,1.037164873491486,1.028582341416547,1.045355939239307),
@(1.033015924296492,1.030643841669302,1.03680649170505,1.028157318198017,1.044899140663327),
@(1.033073708509836,1.030859923977745,1.036996503796713,1.028382599775119,1.045141288443898),
@(1.033298613669594,1.031710756142654,1.037743955312947,1.029270198966425,1.046094810564026),
@(1.033753740601323,1.032698339954248,1.03861002689691,1.030301586055918,1.047201699757723)
)

for ($i = 0; $i -lt $block1.Count; $i++) {
    $row = 2 + $i
    $vals = $block1[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}

for ($i = 0; $i -lt $block2.Count; $i++) {
    $row = 2 + $i
    $vals = $block2[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, 9 + $j).Value = $vals[$j]
    }
}
